$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 2
Set-TextValue $ws.Range("D2") '68.035.65'
$ws.Range("E2").Value = '  +0.10%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.253.79'
$ws.Range("E3").Value = '  -0.78%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '584.24'
$ws.Range("E5").Value = '  +0.47%  '

# Row 6
Set-TextValue $ws.Range("D6") '181.29'
$ws.Range("E6").Value = '  -0.69%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.595'
$ws.Range("E8").Value = '  -1.11%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.134'
$ws.Range("E9").Value = '  -0.04%  '

# Row 10
$ws.Range("E10").Value = '  -1.69%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.422'
$ws.Range("E11").Value = '  +0.66%  '

# Row 12
Set-TextValue $ws.Range("D12") '3.816.17'
$ws.Range("E12").Value = '  -0.86%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.138'
$ws.Range("E13").Value = '  +0.02%  '

# Row 14
Set-TextValue $ws.Range("D14") '28.19'
$ws.Range("E14").Value = '  -1.66%  '

# Row 15
Set-TextValue $ws.Range("D15") '68.085.84'
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("E16").Value = '  +0.54%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.238.54'
$ws.Range("E17").Value = '  -1.30%  '

# Row 18
Set-TextValue $ws.Range("D18") '5.83'
$ws.Range("E18").Value = '  -0.50%  '

# Row 19
Set-TextValue $ws.Range("D19") '13.50'
$ws.Range("E19").Value = '  -0.79%  '

# Row 20
Set-TextValue $ws.Range("D20") '392.56'

# Row 21
Set-TextValue $ws.Range("D21") '7.67'
$ws.Range("E21").Value = '  +0.00%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D22") '1.00'
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D23") '71.41'
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
Set-TextValue $ws.Range("D24") '0.516'
$ws.Range("E24").Value = '  +0.59%  '

# Row 25
Set-TextValue $ws.Range("D25") '0.0000119'
$ws.Range("E25").Value = '  -1.27%  '

# Row 26
$ws.Range("E26").Value = '  +4.20%  '

# Row 27
Set-TextValue $ws.Range("D27") '9.60'
$ws.Range("E27").Value = '  -0.82%  '

# Row 28
Set-TextValue $ws.Range("D28") '0.999'
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.99'
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
Set-TextValue $ws.Range("D30") '5.68'
$ws.Range("E30").Value = '  -0.63%  '

# Row 31
Set-TextValue $ws.Range("D31") '23.03'
$ws.Range("E31").Value = '  +0.90%  '

# Row 32
Set-TextValue $ws.Range("D32") '7.11'
$ws.Range("E32").Value = '  +2.09%  '

# Row 33
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("E34").Value = '  -0.82%  '

# Row 35
Set-TextValue $ws.Range("D35") '164.84'
$ws.Range("E35").Value = '  +0.76%  '

# Row 36
Set-TextValue $ws.Range("D36") '1.49'
$ws.Range("E36").Value = '  -1.53%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.91'
$ws.Range("E37").Value = '  +2.86%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.822'
$ws.Range("E38").Value = '  -3.85%  '

# Row 39
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D39") '26.43'
$ws.Range("E39").Value = '  -1.60%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D40") '4.60'
$ws.Range("E40").Value = '  -1.70%  '

# Row 41
Set-TextValue $ws.Range("D41") '6.55'
$ws.Range("E41").Value = '  -4.50%  '

# Row 42
Set-TextValue $ws.Range("D42") '41.39'
$ws.Range("E42").Value = '  +1.08%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.0689'
$ws.Range("E43").Value = '  +0.66%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.47'
$ws.Range("E44").Value = '  -6.36%  '

# Row 45
Set-TextValue $ws.Range("D45") '342.84'
$ws.Range("E45").Value = '  -3.63%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.595.87'
$ws.Range("E46").Value = '  -4.29%  '

# Row 47
Set-TextValue $ws.Range("D47") '24.72'
$ws.Range("E47").Value = '  -3.53%  '

# Row 48
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
Set-TextValue $ws.Range("D49") '31.81'
$ws.Range("E49").Value = '  +0.82%  '

# Row 50
Set-TextValue $ws.Range("D50") '6.30'
$ws.Range("E50").Value = '  +1.88%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.102'
$ws.Range("E51").Value = '  -1.05%  '
